# Estadisticos Segundo Parcial 26 Mayo
# The "Rescatables" sheet lists rescatable students; reorder the three
# student blocks (CASTRO, RAMOS, CASTILLO) so RAMOS comes first, then
# CASTILLO, then CASTRO -- and update CASTRO's "Reprobadas" count from 3 to 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# New row 5: student 23330051920018 RAMOS / UTRERA / CARLOS DAVID - Control Electronico
$ws.Cells.Item(5, 1).Value = 23330051920018
$ws.Cells.Item(5, 2).Value = "RAMOS"
$ws.Cells.Item(5, 3).Value = "UTRERA"
$ws.Cells.Item(5, 4).Value = "CARLOS DAVID"
$ws.Cells.Item(5, 5).Value = "MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO"
$ws.Cells.Item(5, 6).Value = "4AEM"
$ws.Cells.Item(5, 7).Value = 3

# New row 6: student 23330051920018 RAMOS / UTRERA / CARLOS DAVID - PLC
$ws.Cells.Item(6, 1).Value = 23330051920018
$ws.Cells.Item(6, 2).Value = "RAMOS"
$ws.Cells.Item(6, 3).Value = "UTRERA"
$ws.Cells.Item(6, 4).Value = "CARLOS DAVID"
$ws.Cells.Item(6, 5).Value = "PROGRAMA Y CONECTA CONTROLADORES LÓGICOS PROGRAMABLES (PLC´S)"
$ws.Cells.Item(6, 6).Value = "4AEM"
$ws.Cells.Item(6, 7).Value = 3

# New row 7: student 24330051920340 CASTILLO / GONZALEZ / ANGEL ALBERTO - Pensamiento matematico
$ws.Cells.Item(7, 1).Value = 24330051920340
$ws.Cells.Item(7, 2).Value = "CASTILLO"
$ws.Cells.Item(7, 3).Value = "GONZALEZ"
$ws.Cells.Item(7, 4).Value = "ANGEL ALBERTO"
$ws.Cells.Item(7, 5).Value = "Pensamiento matemático II"
$ws.Cells.Item(7, 6).Value = "2BEM"
$ws.Cells.Item(7, 7).Value = 2

# New row 8: student 23330051920005 CASTRO / ARIAS / OMAR DAVID - Control Electronico
$ws.Cells.Item(8, 1).Value = 23330051920005
$ws.Cells.Item(8, 2).Value = "CASTRO"
$ws.Cells.Item(8, 3).Value = "ARIAS"
$ws.Cells.Item(8, 4).Value = "OMAR DAVID"
$ws.Cells.Item(8, 5).Value = "MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO"
$ws.Cells.Item(8, 6).Value = "4AEM"
$ws.Cells.Item(8, 7).Value = 2

# New row 9: student 23330051920005 CASTRO / ARIAS / OMAR DAVID - PLC
$ws.Cells.Item(9, 1).Value = 23330051920005
$ws.Cells.Item(9, 2).Value = "CASTRO"
$ws.Cells.Item(9, 3).Value = "ARIAS"
$ws.Cells.Item(9, 4).Value = "OMAR DAVID"
$ws.Cells.Item(9, 5).Value = "PROGRAMA Y CONECTA CONTROLADORES LÓGICOS PROGRAMABLES (PLC´S)"
$ws.Cells.Item(9, 6).Value = "4AEM"
$ws.Cells.Item(9, 7).Value = 2
